# qc.vw_ic_reference_all_towns_oby_detail.xlsx
#
# Commit: "Make sheet name configurable in `export_models` and use it to
# tweak sheet name for IC reference files"
#
# The only content-level change the commit describes is the worksheet's
# display name (Sheet1 -> "Query Results"); the rest of the diff
# (absPath, revisionPtr, bookViews window geometry, the stale A3:XFD91500
# selection) is just Excel re-stamping its own last-interactive-session
# chrome on save, not something the export code controls. We still reset
# the selection back to the top-left cell so the sheet doesn't keep
# carrying the old, now-meaningless A3:XFD91500 selection forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Query Results).
$ws.Name = "Query Results"

# Drop the leftover selection (A3:XFD91500) by reselecting the home cell.
$ws.Activate()
$ws.Range("A1").Select()
